$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.84444160904917
$ws.Range("C2").Value = 7.404700331357763
$ws.Range("D2").Value = 13.79735813243375
$ws.Range("E2").Value = 14.30943522161497
$ws.Range("G2").Value = 38.75699452027821
$ws.Range("H2").Value = 16.92050874979337
$ws.Range("I2").Value = 26.69889457444484
$ws.Range("J2").Value = 8.574804373849281
$ws.Range("K2").Value = 10.75487860213764
$ws.Range("L2").Value = 12.1201095686154
$ws.Range("N2").Value = 20.30674022092482
$ws.Range("O2").Value = 27.06042853053681
$ws.Range("B3").Value = 13.64082652023611
$ws.Range("C3").Value = 7.369289078718641
$ws.Range("D3").Value = 13.78961947914351
$ws.Range("E3").Value = 14.32996176696378
$ws.Range("G3").Value = 38.85317870845344
$ws.Range("H3").Value = 16.96609453911172
$ws.Range("I3").Value = 26.78430182364682
$ws.Range("J3").Value = 8.585209836460187
$ws.Range("K3").Value = 10.60936204762067
$ws.Range("L3").Value = 12.12088230673
$ws.Range("N3").Value = 20.36238588133515
$ws.Range("O3").Value = 27.13680343401371
$ws.Range("B4").Value = 13.51646607814988
$ws.Range("C4").Value = 7.347232337629229
$ws.Range("D4").Value = 13.78729597499428
$ws.Range("E4").Value = 14.34445389550186
$ws.Range("G4").Value = 38.92095896967271
$ws.Range("H4").Value = 16.99626509354551
$ws.Range("I4").Value = 26.84074884013348
$ws.Range("J4").Value = 8.591961867301915
$ws.Range("K4").Value = 10.52051430689779
$ws.Range("L4").Value = 12.12294064291302
$ws.Range("N4").Value = 20.39818425652604
$ws.Range("O4").Value = 27.18814753736513
$ws.Range("B5").Value = 13.46601466243442
$ws.Range("C5").Value = 7.338167656054208
$ws.Range("D5").Value = 13.78696189991491
$ws.Range("E5").Value = 14.35083500861138
$ws.Range("G5").Value = 38.95076917523129
$ws.Range("H5").Value = 17.00910861204253
$ws.Range("I5").Value = 26.86475923195768
$ws.Range("J5").Value = 8.594804923438158
$ws.Range("K5").Value = 10.48447421182713
$ws.Range("L5").Value = 12.12417894351116
$ws.Range("N5").Value = 20.41318385532617
$ws.Range("O5").Value = 27.21018912301982
$ws.Range("B6").Value = 13.45765269100808
$ws.Range("C6").Value = 7.33665793512486
$ws.Range("D6").Value = 13.7869434954776
$ws.Range("E6").Value = 14.35192331777392
$ws.Range("G6").Value = 38.95585122826004
$ws.Range("H6").Value = 17.01127442535152
$ws.Range("I6").Value = 26.86880701462886
$ws.Range("J6").Value = 8.595282547739822
$ws.Range("K6").Value = 10.47850098994406
$ws.Range("L6").Value = 12.12440872574525
$ws.Range("N6").Value = 20.41569941329748
$ws.Range("O6").Value = 27.21391664929073
$ws.Range("B7").Value = 13.51578467574704
$ws.Range("C7").Value = 7.347110393998101
$ws.Range("D7").Value = 13.78728898575063
$ws.Range("E7").Value = 14.34453802768678
$ws.Range("G7").Value = 38.92135214242224
$ws.Range("H7").Value = 16.99643608321291
$ws.Range("I7").Value = 26.84106857186127
$ws.Range("J7").Value = 8.591999838725449
$ws.Range("K7").Value = 10.52002753186843
$ws.Range("L7").Value = 12.12295572391009
$ws.Range("N7").Value = 20.39838487879023
$ws.Range("O7").Value = 27.18844026996464
$ws.Range("B8").Value = 13.7741343767826
$ws.Range("C8").Value = 7.392557278822136
$ws.Range("D8").Value = 13.79418732127658
$ws.Range("E8").Value = 14.31612108633412
$ws.Range("G8").Value = 38.78834599827315
$ws.Range("H8").Value = 16.9357743797251
$ws.Range("I8").Value = 26.72751158830193
$ws.Range("J8").Value = 8.578316986553101
$ws.Range("K8").Value = 10.70462492516627
$ws.Range("L8").Value = 12.12004804810841
$ws.Range("N8").Value = 20.32558888479397
$ws.Range("O8").Value = 27.08583882884624
$ws.Range("B9").Value = 14.2833533661984
$ws.Range("C9").Value = 7.479073666834096
$ws.Range("D9").Value = 13.82687099076366
$ws.Range("E9").Value = 14.27535823746094
$ws.Range("G9").Value = 38.59690608782907
$ws.Range("H9").Value = 16.83410103845438
$ws.Range("I9").Value = 26.53660484222059
$ws.Range("J9").Value = 8.554353625956114
$ws.Range("K9").Value = 11.06881961050743
$ws.Range("L9").Value = 12.12685570872672
$ws.Range("N9").Value = 20.19572920760501
$ws.Range("O9").Value = 26.91996213927392
$ws.Range("B10").Value = 14.65557586125284
$ws.Range("C10").Value = 7.54090067634712
$ws.Range("D10").Value = 13.86239705174102
$ws.Range("E10").Value = 14.25449781758372
$ws.Range("G10").Value = 38.49876809522833
$ws.Range("H10").Value = 16.76991269099793
$ws.Range("I10").Value = 26.41570091920739
$ws.Range("J10").Value = 8.538480052377096
$ws.Range("K10").Value = 11.3353928992502
$ws.Range("L10").Value = 12.13940087062166
$ws.Range("N10").Value = 20.10810463750351
$ws.Range("O10").Value = 26.81965142767297
$ws.Range("B11").Value = 14.82378326483861
$ws.Range("C11").Value = 7.568618747082582
$ws.Range("D11").Value = 13.88101706912943
$ws.Range("E11").Value = 14.24697301383739
$ws.Range("G11").Value = 38.46339025499002
$ws.Range("H11").Value = 16.74298870100199
$ws.Range("I11").Value = 26.36489789300911
$ws.Range("J11").Value = 8.531631405030872
$ws.Range("K11").Value = 11.45596348762468
$ws.Range("L11").Value = 12.14672599141941
$ws.Range("N11").Value = 20.06991570544049
$ws.Range("O11").Value = 26.77870317184563
$ws.Range("B12").Value = 14.8872609671805
$ws.Range("C12").Value = 7.579053555093217
$ws.Range("D12").Value = 13.88841749848507
$ws.Range("E12").Value = 14.24440528283666
$ws.Range("G12").Value = 38.45132808139947
$ws.Range("H12").Value = 16.73312011917444
$ws.Range("I12").Value = 26.34626341549147
$ws.Range("J12").Value = 8.529091274308628
$ws.Range("K12").Value = 11.50148145014015
$ws.Range("L12").Value = 12.14973059869579
$ws.Range("N12").Value = 20.05569376981167
$ws.Range("O12").Value = 26.76387093998426
$ws.Range("B13").Value = 14.87360053882162
$ws.Range("C13").Value = 7.57680900956991
$ws.Range("D13").Value = 13.88680821237505
$ws.Range("E13").Value = 14.24494577180414
$ws.Range("G13").Value = 38.45386649733521
$ws.Range("H13").Value = 16.7352309581869
$ws.Range("I13").Value = 26.35024984106457
$ws.Range("J13").Value = 8.52963596995995
$ws.Range("K13").Value = 11.49168516337395
$ws.Range("L13").Value = 12.14907327298258
$ws.Range("N13").Value = 20.05874608782738
$ws.Range("O13").Value = 26.76703534373278
$ws.Range("B14").Value = 14.82901032281749
$ws.Range("C14").Value = 7.569478472320593
$ws.Range("D14").Value = 13.88161892880161
$ws.Range("E14").Value = 14.24675612305644
$ws.Range("G14").Value = 38.46237113338755
$ws.Range("H14").Value = 16.74217025590574
$ws.Range("I14").Value = 26.36335272755607
$ws.Range("J14").Value = 8.531421360011704
$ws.Range("K14").Value = 11.45971129576592
$ws.Range("L14").Value = 12.14696856951965
$ws.Range("N14").Value = 20.06874086652051
$ws.Range("O14").Value = 26.77746940497311
$ws.Range("B15").Value = 14.80166735152226
$ws.Range("C15").Value = 7.564980223034524
$ws.Range("D15").Value = 13.8784857198822
$ws.Range("E15").Value = 14.24790168342848
$ws.Range("G15").Value = 38.4677543367532
$ws.Range("H15").Value = 16.7464633461667
$ws.Range("I15").Value = 26.37145722186678
$ws.Range("J15").Value = 8.532521897681496
$ws.Range("K15").Value = 11.44010705848785
$ws.Range("L15").Value = 12.14570936705329
$ws.Range("N15").Value = 20.07489410136847
$ws.Range("O15").Value = 26.7839483548618
$ws.Range("B16").Value = 14.64455560285597
$ws.Range("C16").Value = 7.539080794745439
$ws.Range("D16").Value = 13.86122932253284
$ws.Range("E16").Value = 14.25502905147534
$ws.Range("G16").Value = 38.5012667740117
$ws.Range("H16").Value = 16.77171800414013
$ws.Range("I16").Value = 26.41910546916473
$ws.Range("J16").Value = 8.538935099421559
$ws.Range("K16").Value = 11.32749593271604
$ws.Range("L16").Value = 12.13895455559345
$ws.Range("N16").Value = 20.11063392855893
$ws.Range("O16").Value = 26.82242178149315
$ws.Range("B17").Value = 14.54784426847662
$ws.Range("C17").Value = 7.523085974371021
$ws.Range("D17").Value = 13.85127003325787
$ws.Range("E17").Value = 14.25990413945039
$ws.Range("G17").Value = 38.52420049783387
$ws.Range("H17").Value = 16.78779355735632
$ws.Range("I17").Value = 26.4494110436429
$ws.Range("J17").Value = 8.542964582746317
$ws.Range("K17").Value = 11.25820611993152
$ws.Range("L17").Value = 12.13522381857738
$ws.Range("N17").Value = 20.13298663384076
$ws.Range("O17").Value = 26.8472239293543
$ws.Range("B18").Value = 14.49211653153036
$ws.Range("C18").Value = 7.513848005352338
$ws.Range("D18").Value = 13.8457734399121
$ws.Range("E18").Value = 14.26289313888241
$ws.Range("G18").Value = 38.53826339223729
$ws.Range("H18").Value = 16.79725399848178
$ws.Range("I18").Value = 26.46723702643428
$ws.Range("J18").Value = 8.545317293552959
$ws.Range("K18").Value = 11.21828933519962
$ws.Range("L18").Value = 12.13323044138762
$ws.Range("N18").Value = 20.14600074940216
$ws.Range("O18").Value = 26.86193028486444
$ws.Range("B19").Value = 14.47323237009942
$ws.Range("C19").Value = 7.51071371871115
$ws.Range("D19").Value = 13.84395230891966
$ws.Range("E19").Value = 14.26393695534477
$ws.Range("G19").Value = 38.54317453971107
$ws.Range("H19").Value = 16.80049394031553
$ws.Range("I19").Value = 26.47334044164097
$ws.Range("J19").Value = 8.54611990957363
$ws.Range("K19").Value = 11.20476458331829
$ws.Range("L19").Value = 12.13258175821593
$ws.Range("N19").Value = 20.15043417596217
$ws.Range("O19").Value = 26.86698529885975
$ws.Range("B20").Value = 14.55815033158476
$ws.Range("C20").Value = 7.524792621560275
$ws.Range("D20").Value = 13.85230626299543
$ws.Range("E20").Value = 14.25936603890987
$ws.Range("G20").Value = 38.52166889456194
$ws.Range("H20").Value = 16.78606011957202
$ws.Range("I20").Value = 26.44614408038983
$ws.Range("J20").Value = 8.542532010653032
$ws.Range("K20").Value = 11.26558896178243
$ws.Range("L20").Value = 12.13560519832022
$ws.Range("N20").Value = 20.13059086425748
$ws.Range("O20").Value = 26.84453807182654
$ws.Range("B21").Value = 14.84211393926155
$ws.Range("C21").Value = 7.571633318479325
$ws.Range("D21").Value = 13.88313369852821
$ws.Range("E21").Value = 14.24621673891361
$ws.Range("G21").Value = 38.45983687837677
$ws.Range("H21").Value = 16.7401231448996
$ws.Range("I21").Value = 26.35948771235659
$ws.Range("J21").Value = 8.53089550290578
$ws.Range("K21").Value = 11.46910689084534
$ws.Range("L21").Value = 12.14758052590629
$ws.Range("N21").Value = 20.065798669131
$ws.Range("O21").Value = 26.7743863729295
$ws.Range("B22").Value = 15.02639997135537
$ws.Range("C22").Value = 7.601887083218731
$ws.Range("D22").Value = 13.90531580981723
$ws.Range("E22").Value = 14.23926475432081
$ws.Range("G22").Value = 38.42720608467462
$ws.Range("H22").Value = 16.71200625377251
$ws.Range("I22").Value = 26.30637059630907
$ws.Range("J22").Value = 8.523600957430485
$ws.Range("K22").Value = 11.60128706908592
$ws.Range("L22").Value = 12.15675104246611
$ws.Range("N22").Value = 20.02484806485141
$ws.Range("O22").Value = 26.73246686042618
$ws.Range("B23").Value = 14.92818063356861
$ws.Range("C23").Value = 7.585773897281511
$ws.Range("D23").Value = 13.89329209522975
$ws.Range("E23").Value = 14.24282520294359
$ws.Range("G23").Value = 38.44390927711618
$ws.Range("H23").Value = 16.72683850029776
$ws.Range("I23").Value = 26.33439830660813
$ws.Range("J23").Value = 8.527465852466205
$ws.Range("K23").Value = 11.53082875304056
$ws.Range("L23").Value = 12.15173427454591
$ws.Range("N23").Value = 20.0465768846531
$ws.Range("O23").Value = 26.75448047287837
$ws.Range("B24").Value = 14.55349135190966
$ws.Range("C24").Value = 7.524021178020155
$ws.Range("D24").Value = 13.85183706949984
$ws.Range("E24").Value = 14.25960873397653
$ws.Range("G24").Value = 38.52281069770467
$ws.Range("H24").Value = 16.78684312649543
$ws.Range("I24").Value = 26.44761981945853
$ws.Range("J24").Value = 8.542727464009666
$ws.Range("K24").Value = 11.26225142870941
$ws.Range("L24").Value = 12.13543230457303
$ws.Range("N24").Value = 20.13167348309632
$ws.Range("O24").Value = 26.84575095483457
$ws.Range("B25").Value = 14.1456822548684
$ws.Range("C25").Value = 7.455962873918215
$ws.Range("D25").Value = 13.81599409739527
$ws.Range("E25").Value = 14.28478704747724
$ws.Range("G25").Value = 38.64124357357971
$ws.Range("H25").Value = 16.85975866208565
$ws.Range("I25").Value = 26.58484994936596
$ws.Range("J25").Value = 8.560530953707785
$ws.Range("K25").Value = 10.97030386312969
$ws.Range("L25").Value = 12.12368285300153
$ws.Range("N25").Value = 20.22948736175213
$ws.Range("O25").Value = 26.96105177084742
